# Integrate retry logic and data provider enhancements; add method
# interceptor and annotation transformer for test execution control.
#
# Functional edit: add a new "TestsRunner" worksheet (after "Sheet1") that
# drives which TestNG test cases execute, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# --- add the new sheet and put it after Sheet1 in tab order -----------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "TestsRunner"

# fetch live handles (post-Add) before reordering, then move into place
$sheet1 = $wb.Worksheets.Item("Sheet1")
$newSheet2 = $wb.Worksheets.Item("TestsRunner")
$newSheet2.Move($null, $sheet1)

# re-fetch a live handle to the moved sheet before mutating it
$ws = $wb.Worksheets.Item("TestsRunner")

# --- header row --------------------------------------------------------------
$ws.Range("A1").Value = "TestCase"
$ws.Range("B1").Value = "Description"

# --- test case rows -----------------------------------------------------------
$ws.Range("A2").Value = "loginLogoutTest"
$ws.Range("B2").Value = "validate OrangeHRM login and logout functionality"

$ws.Range("C3").Value = "yes"
$ws.Range("A3").Value = "homePageTitleTest"
$ws.Range("C2").Value = "no"
$ws.Range("B3").Value = "validate title of home page"
$ws.Range("C1").Value = "Execute"

# --- column widths -------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.166666666666668
$ws.Columns.Item(2).ColumnWidth = 47.5

# --- make TestsRunner the active/selected sheet, with C4 selected ------------
$ws.Activate()
$ws.Range("C4").Select() | Out-Null
